$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

$ws.Range("C30").Copy($ws.Range("D14"))
$ws.Range("E30").Copy($ws.Range("E14"))
$ws.Range("G14").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 2
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = 14.285714285714
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 166.666666666667
$ws.Range("N15").Value = -11.111111111111
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 7.142857142857
$ws.Range("I16").Value = 78
$ws.Range("J16").Value = 73
$ws.Range("K16").Value = 6.849315068493
$ws.Range("L16").Value = 59.183673469387
$ws.Range("M16").Value = -7.142857142857
$ws.Range("N16").Value = -74.426229508196
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 21.428571428571
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = 30.232558139534
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 23.762376237623
$ws.Range("L17").Value = 42.045454545454
$ws.Range("M17").Value = 52.439024390243
$ws.Range("N17").Value = -2.34375
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = -2.439024390243
$ws.Range("L18").Value = 37.931034482758
$ws.Range("M18").Value = -47.368421052631
$ws.Range("N18").Value = -90.543735224586
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -12.280701754386
$ws.Range("I19").Value = 118
$ws.Range("J19").Value = 135
$ws.Range("K19").Value = -12.592592592592
$ws.Range("L19").Value = 11.320754716981
$ws.Range("M19").Value = 12.380952380952
$ws.Range("N19").Value = -24.840764331210
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 22
$ws.Range("H20").Value = 4.761904761904
$ws.Range("I20").Value = 69
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = 32.692307692307
$ws.Range("L20").Value = 115.625
$ws.Range("M20").Value = 122.58064516129
$ws.Range("N20").Value = -77.450980392156
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 9.756097560975
$ws.Range("G21").Value = 173
$ws.Range("H21").Value = 0.578034682080
$ws.Range("I21").Value = 438
$ws.Range("J21").Value = 414
$ws.Range("K21").Value = 5.797101449275
$ws.Range("L21").Value = 41.290322580645
$ws.Range("M21").Value = 14.659685863874
$ws.Range("N21").Value = -67.264573991031
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 2
$ws.Range("C30").Copy($ws.Range("D22"))
$ws.Range("E30").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 133.333333333333
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = 12.5
$ws.Range("L22").Value = 80
$ws.Range("M22").Value = 200
$ws.Range("C30").Copy($ws.Range("D23"))
$ws.Range("E30").Copy($ws.Range("E23"))
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = -32
$ws.Range("F24").Value = 186
$ws.Range("G24").Value = 234
$ws.Range("H24").Value = -20.512820512820
$ws.Range("I24").Value = 434
$ws.Range("J24").Value = 622
$ws.Range("K24").Value = -30.225080385852
$ws.Range("L24").Value = 78.600823045267
$ws.Range("M24").Value = 76.422764227642
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -7.142857142857
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -1.886792452830
$ws.Range("I25").Value = 133
$ws.Range("J25").Value = 135
$ws.Range("K25").Value = -1.481481481481
$ws.Range("L25").Value = 26.666666666666
$ws.Range("M25").Value = -13.071895424836
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 2
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 2
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -44.444444444444
$ws.Range("L26").Value = 42.857142857142
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = 63.636363636363
$ws.Range("L27").Value = 5.882352941176
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = 25
$ws.Range("M28").Value = 66.666666666666
$ws.Range("N28").Value = -61.538461538461
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -44.444444444444
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = 66.666666666666
$ws.Range("N29").Value = -44.444444444444
